$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header/count values changed ---
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# --- Row 2 (CON): D2 value removed, new values placed in B2 and C2 ---
$ws.Range("D2").ClearContents() | Out-Null
$ws.Range("B2").Value = 41.672936990180958
$ws.Range("C2").Value = 21.834834647764225

# --- Row 3 (STR): B3 value removed, C3 value replaced ---
$ws.Range("B3").ClearContents() | Out-Null
$ws.Range("C3").Value = 40.595687312940726

# --- Selection narrowed from B1:AY3 to B1:E3 ---
$ws.Range("B1:E3").Select() | Out-Null
